# Fix Training Data Issue (#48)
# Data was taken from 1 day off due to way NBA stats were shown.
# The "Date" column (BF) on Sheet1 held values formatted like "5-9-2013-14"
# (an awkward/ambiguous mash-up of game-date + season) for every data row
# (rows 2-31). Rewrite each one as an ISO-style date string "2014-05-09".
#
# Cells are written as literal text (not re-parsed as an Excel date serial)
# by forcing a leading quote-prefix, then immediately reapplying the
# "Normal" cell style so no stray style/formatting is attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldDate = "5-9-2013-14"
$newDate = "2014-05-09"

$firstRow = 2
$lastRow = 31
$col = "BF"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$col$row")
    if ($cell.Value() -eq $oldDate) {
        $cell.Value = "'" + $newDate
        $cell.Style = "Normal"
    }
}
